$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2 and 3 with new survey responses ---
$ws.Cells.Item(2,1).Value = 45242.41506680555
$ws.Cells.Item(2,2).Value = "Não"
$ws.Cells.Item(2,3).Value = "Não"
$ws.Cells.Item(2,4).Value = "Hoje é um dia incrível, cheio de boas energias!"

$ws.Cells.Item(3,1).Value = 45242.41529898148
$ws.Cells.Item(3,2).Value = "Não"
$ws.Cells.Item(3,3).Value = "Sim"
$ws.Cells.Item(3,4).Value = "Não consigo superar a frustração dessa situação"

# --- Add 3 new response rows (4, 5, 6) ---
# Copy formatting from row 2 (an existing, fully-styled response row) so the
# new rows pick up the same cell styles (date format in col A, text style in B:D)
# instead of the engine creating brand-new style entries.

$ws.Range("A2:D2").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$ws.Cells.Item(4,1).Value = 45242.41548430556
$ws.Cells.Item(4,2).Value = "Sim"
$ws.Cells.Item(4,3).Value = "Não"
$ws.Cells.Item(4,4).Value = "Recebi uma notícia maravilhosa que me deixou radiante!"

$ws.Range("A2:D2").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
$ws.Cells.Item(5,1).Value = 45242.41566929399
$ws.Cells.Item(5,2).Value = "Não"
$ws.Cells.Item(5,3).Value = "Sim"
$ws.Cells.Item(5,4).Value = "Sinto-me um pouco ansioso com os próximos desafios"

$ws.Range("A2:D2").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Cells.Item(6,1).Value = 45242.4158944213
$ws.Cells.Item(6,2).Value = "Não"
$ws.Cells.Item(6,3).Value = "Sim"
$ws.Cells.Item(6,4).Value = "A tristeza parece persistir, não sei como lidar com isso"

$excel.CutCopyMode = 0

Write-Host "edit complete"
